# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.377.32'
$ws.Range("E2").Value = '  -1.04%  '

$ws.Range("D3").Value = '2.612.02'
$ws.Range("E3").Value = '  +2.65%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = '''306.43'
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").Value = '''100.07'
$ws.Range("E6").Value = '  -3.48%  '

$ws.Range("D7").Value = '''0.601'
$ws.Range("E7").Value = '  -1.14%  '

$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '''0.579'
$ws.Range("E9").Value = '  +1.60%  '

$ws.Range("D10").Value = '''39.40'
$ws.Range("E10").Value = '  +0.90%  '

$ws.Range("D11").Value = '''54.20'
$ws.Range("E11").Value = '  -0.52%  '

$ws.Range("D12").Value = '''0.0842'
$ws.Range("E12").Value = '  +1.58%  '

$ws.Range("D13").Value = '''8.11'
$ws.Range("E13").Value = '  +2.10%  '

$ws.Range("D14").Value = '3.002.62'
$ws.Range("E14").Value = '  +2.42%  '

$ws.Range("E15").Value = '  +0.75%  '

$ws.Range("D16").Value = '2.606.73'
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").Value = '''0.920'
$ws.Range("E17").Value = '  +2.36%  '

$ws.Range("D18").Value = '''14.95'
$ws.Range("E18").Value = '  -0.53%  '

$ws.Range("D19").Value = '46.459.40'
$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").Value = '''0.0000101'
$ws.Range("E20").Value = '  +1.51%  '

$ws.Range("D21").Value = '''12.96'
$ws.Range("E21").Value = '  -8.53%  '

$ws.Range("D22").Value = '''6.72'
$ws.Range("E22").Value = '  +1.84%  '

$ws.Range("D23").Value = '''71.49'
$ws.Range("E23").Value = '  +2.24%  '

$ws.Range("D24").Value = '''272.91'
$ws.Range("E24").Value = '  +6.89%  '

$ws.Range("D25").Value = '''3.03'
$ws.Range("E25").Value = '  +1.88%  '

$ws.Range("D26").Value = '''2.17'
$ws.Range("E26").Value = '  +2.32%  '

$ws.Range("D27").Value = '''29.52'
$ws.Range("E27").Value = '  +19.93%  '

$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").Value = '''4.02'
$ws.Range("E29").Value = '  -0.68%  '

$ws.Range("D30").Value = '''10.59'
$ws.Range("E30").Value = '  +0.78%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '''38.49'
$ws.Range("E31").Value = '  -8.61%  '

$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = '''2.21'
$ws.Range("E32").Value = '  -2.64%  '

$ws.Range("D33").Value = '''6.33'
$ws.Range("E33").Value = '  +5.50%  '

$ws.Range("D34").Value = '''3.67'
$ws.Range("E34").Value = '  -3.13%  '

$ws.Range("E35").Value = '  -2.69%  '

$ws.Range("D36").Value = '''2.23'
$ws.Range("E36").Value = '  +1.73%  '

$ws.Range("D37").Value = '''0.0836'

$ws.Range("D38").Value = '''151.56'
$ws.Range("E38").Value = '  +0.93%  '

$ws.Range("E39").Value = '  +0.84%  '

$ws.Range("D40").Value = '''0.123'
$ws.Range("E40").Value = '  +1.37%  '

$ws.Range("D41").Value = '''23.31'
$ws.Range("E41").Value = '  +32.55%  '

$ws.Range("D42").Value = '''15.86'
$ws.Range("E42").Value = '  -3.34%  '

$ws.Range("D43").Value = '''0.0329'
$ws.Range("E43").Value = '  +0.94%  '

$ws.Range("D44").Value = '''3.61'
$ws.Range("E44").Value = '  +1.00%  '

$ws.Range("D45").Value = '''4.07'
$ws.Range("E45").Value = '  -5.23%  '

$ws.Range("D46").Value = '2.128.59'
$ws.Range("E46").Value = '  +6.23%  '

$ws.Range("E47").Value = '  +0.15%  '

$ws.Range("D48").Value = '''94.04'
$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("D49").Value = '''9.55'
$ws.Range("E49").Value = '  +7.69%  '

$ws.Range("D50").Value = '''1.78'
$ws.Range("E50").Value = '  -6.38%  '

$ws.Range("D51").Value = '''108.60'
$ws.Range("E51").Value = '  +1.57%  '

